$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster changed from "MuSCs" to "ECs" for all data rows
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"

# Target cluster values (texts unchanged, kept explicit for clarity)
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"

# Row 2 numeric updates
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02551366666666667
$ws.Range("H2").Value = 0.076541
$ws.Range("M2").Value = 1.684496
$ws.Range("N2").Value = 5.053488
$ws.Range("O2").Value = 0.6423607101334534
$ws.Range("P2").Value = 0.7282461611889918
$ws.Range("Q2").Value = 0.04297766944533334
$ws.Range("R2").Value = 0.386799025008
$ws.Range("S2").Value = 0.6423607101334534
$ws.Range("T2").Value = 0.7282461611889918

# Row 3 numeric updates
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02551366666666667
$ws.Range("H3").Value = 0.076541
$ws.Range("M3").Value = 0.9277985
$ws.Range("N3").Value = 1.855597
$ws.Range("O3").Value = 0.3538039290807178
$ws.Range("P3").Value = 0.2674056793968462
$ws.Range("Q3").Value = 0.02367154166283333
$ws.Range("R3").Value = 0.142029249977
$ws.Range("S3").Value = 0.3538039290807178
$ws.Range("T3").Value = 0.2674056793968462

# Row 4 numeric updates
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02551366666666667
$ws.Range("H4").Value = 0.076541
$ws.Range("M4").Value = 0.01005766666666667
$ws.Range("N4").Value = 0.030173
$ws.Range("O4").Value = 0.003835360785828855
$ws.Range("P4").Value = 0.004348159414162149
$ws.Range("Q4").Value = 0.0002566079547777778
$ws.Range("R4").Value = 0.002309471593
$ws.Range("S4").Value = 0.003835360785828855
$ws.Range("T4").Value = 0.004348159414162149
